$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.603177
$ws.Range("H2").Value = 19.809531
$ws.Range("I2").Value = 0.5135477412645301
$ws.Range("J2").Value = 0.5135477412645302
$ws.Range("M2").Value = 0.06624833333333334
$ws.Range("N2").Value = 0.198745
$ws.Range("Q2").Value = 0.437449470955
$ws.Range("R2").Value = 3.937045238595
$ws.Range("S2").Value = 0.5135477412645301
$ws.Range("T2").Value = 0.5135477412645302

# Row 3
$ws.Range("I3").Value = 0.02944398858046029
$ws.Range("J3").Value = 0.0294439885804603
$ws.Range("M3").Value = 0.06624833333333334
$ws.Range("N3").Value = 0.198745
$ws.Range("Q3").Value = 0.02508093443388889
$ws.Range("R3").Value = 0.225728409905
$ws.Range("S3").Value = 0.02944398858046029
$ws.Range("T3").Value = 0.0294439885804603

# Row 4
$ws.Range("G4").Value = 3.441487333333333
$ws.Range("H4").Value = 10.324462
$ws.Range("I4").Value = 0.2676541983690312
$ws.Range("J4").Value = 0.2676541983690313
$ws.Range("M4").Value = 0.06624833333333334
$ws.Range("N4").Value = 0.198745
$ws.Range("Q4").Value = 0.2279928000211111
$ws.Range("R4").Value = 2.05193520019
$ws.Range("S4").Value = 0.2676541983690312
$ws.Range("T4").Value = 0.2676541983690313

# Row 5
$ws.Range("G5").Value = 2.434707333333333
$ws.Range("H5").Value = 7.304122
$ws.Range("I5").Value = 0.1893540717859783
$ws.Range("J5").Value = 0.1893540717859783
$ws.Range("M5").Value = 0.06624833333333334
$ws.Range("N5").Value = 0.198745
$ws.Range("Q5").Value = 0.1612953029877778
$ws.Range("R5").Value = 1.45165772689
$ws.Range("S5").Value = 0.1893540717859783
$ws.Range("T5").Value = 0.1893540717859783
